# INDUSTOWER.NS.xlsx weekly stock data — "break out stock.yaml completed"
#
# 1) Fix three pre-existing cells that were mis-populated (Q58, Q64, O599)
#    and two trailing rows whose "backup" column was left as an empty
#    inline string instead of a numeric 0 (R601, R602).
# 2) Append 13 new weekly candles (rows 603-615, 2024-07-01 .. 2024-09-23).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Point fixes on existing rows ---------------------------------------
$ws.Cells.Item(58, 17).Value = 0     # Q58  detect_structure: 2 -> 0
$ws.Cells.Item(64, 17).Value = 0     # Q64  detect_structure: 1 -> 0
$ws.Cells.Item(599, 15).Value = 2    # O599 isPivot: 0 -> 2
$ws.Cells.Item(601, 18).Value = 0    # R601 backup: "" -> 0
$ws.Cells.Item(602, 18).Value = 0    # R602 backup: "" -> 0

# --- 2) Append new weekly rows 603-615 --------------------------------------
# Columns: Row, Datetime(serial), Open, High, Low, Close, AdjClose, Volume,
#          Year, Month, Day, Hour, Minute, Second, Week, isPivot,
#          two_line_structure, detect_structure
$newRows = @(
    @(603, 45474, 377, 408.3999938964844, 372.75, 395.6499938964844, 395.6499938964844, 122828989, 2024, 7, 1, 0, 0, 0, 27, 0, 0, 0),
    @(604, 45481, 396, 398.8500061035156, 373.5499877929688, 392, 392, 70078542, 2024, 7, 8, 0, 0, 0, 28, 0, 0, 0),
    @(605, 45488, 392.0499877929688, 424.7000122070312, 388, 409.5, 409.5, 84567404, 2024, 7, 15, 0, 0, 0, 29, 0, 0, 1),
    @(606, 45495, 409.6000061035156, 447.3500061035156, 403.6499938964844, 444.8500061035156, 444.8500061035156, 153335008, 2024, 7, 22, 0, 0, 0, 30, 0, 0, 0),
    @(607, 45502, 448.8999938964844, 453.2999877929688, 420.1000061035156, 421.1000061035156, 421.1000061035156, 62697383, 2024, 7, 29, 0, 0, 0, 31, 0, 0, 0),
    @(608, 45509, 417.8999938964844, 427.2999877929688, 407, 414.9500122070312, 414.9500122070312, 60777723, 2024, 8, 5, 0, 0, 0, 32, 0, 0, 0),
    @(609, 45516, 415, 421, 400.6499938964844, 411.7999877929688, 411.7999877929688, 22729954, 2024, 8, 12, 0, 0, 0, 33, 0, 0, 0),
    @(610, 45523, 413.7000122070312, 439.2000122070312, 412.1000061035156, 434.3500061035156, 434.3500061035156, 38701295, 2024, 8, 19, 0, 0, 0, 34, 0, 0, 0),
    @(611, 45530, 435.4500122070312, 460.3500061035156, 428.5499877929688, 458.5, 458.5, 61526843, 2024, 8, 26, 0, 0, 0, 35, 0, 0, 0),
    @(612, 45537, 460, 460, 414, 423.1000061035156, 423.1000061035156, 56994489, 2024, 9, 2, 0, 0, 0, 36, 0, 0, 0),
    @(613, 45544, 420.9500122070312, 438.6499938964844, 416.7999877929688, 428.4500122070312, 428.4500122070312, 38287956, 2024, 9, 9, 0, 0, 0, 37, 0, 0, 0),
    @(614, 45551, 428.4500122070312, 433.6000061035156, 366.3500061035156, 388.25, 388.25, 141374472, 2024, 9, 16, 0, 0, 0, 38, 0, 0, 0),
    @(615, 45558, 399.8999938964844, 412.6000061035156, 385.0499877929688, 392.3999938964844, 392.3999938964844, 67509516, 2024, 9, 23, 0, 0, 0, 39, 0, 0, 0)
)

foreach ($row in $newRows) {
    $r = $row[0]

    # Column A carries the same custom datetime number format as every
    # other row in the sheet.
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($r, 2).Value = $row[2]    # Open
    $ws.Cells.Item($r, 3).Value = $row[3]    # High
    $ws.Cells.Item($r, 4).Value = $row[4]    # Low
    $ws.Cells.Item($r, 5).Value = $row[5]    # Close
    $ws.Cells.Item($r, 6).Value = $row[6]    # Adj Close
    $ws.Cells.Item($r, 7).Value = $row[7]    # Volume
    $ws.Cells.Item($r, 8).Value = $row[8]    # Year
    $ws.Cells.Item($r, 9).Value = $row[9]    # Month
    $ws.Cells.Item($r, 10).Value = $row[10]  # Day
    $ws.Cells.Item($r, 11).Value = $row[11]  # Hour
    $ws.Cells.Item($r, 12).Value = $row[12]  # Minute
    $ws.Cells.Item($r, 13).Value = $row[13]  # Second
    $ws.Cells.Item($r, 14).Value = $row[14]  # Week
    $ws.Cells.Item($r, 15).Value = $row[15]  # isPivot
    $ws.Cells.Item($r, 16).Value = $row[16]  # two_line_structure
    $ws.Cells.Item($r, 17).Value = $row[17]  # detect_structure
    # Column R ("backup") is intentionally left blank for these new rows,
    # matching the source diff (an empty inline string, same as the rest
    # of the freshly scraped rows before they get back-filled later).
}
